# Update attendance/view counts on the "展览" (Exhibition) and
# "全部类型" (All types) worksheets to match the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F17").Value = 5709
$wsExhibit.Range("F19").Value = 236
$wsExhibit.Range("F20").Value = 1350
$wsExhibit.Range("F22").Value = 351

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F18").Value = 5709
$wsAll.Range("F21").Value = 236
$wsAll.Range("F22").Value = 1350
$wsAll.Range("F24").Value = 351
